$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rewrite the sample-mapping table rows 32-36 --------------------------
# The "getFirstOrNullObject" sample (Field.getFirstOrNullObject) is being
# remapped from the Field class to the FieldCollection class, so it moves
# down to sit next to the other FieldCollection row ("items"), and the
# remaining Field rows (code/parentBody/result) shift up to fill rows 32-34.

$ws.Range("A32").Value = "Field"
$ws.Range("B32").Value = "code"
$ws.Range("C32").ClearContents()
$ws.Range("D32").Value = "word-manage-fields"
$ws.Range("E32").Value = "getFirstField"

$ws.Range("A33").Value = "Field"
$ws.Range("B33").Value = "parentBody"
$ws.Range("C33").ClearContents()
$ws.Range("D33").Value = "word-manage-fields"
$ws.Range("E33").Value = "getParentBodyOfFirstField"

$ws.Range("A34").Value = "Field"
$ws.Range("B34").Value = "result"
$ws.Range("C34").ClearContents()
$ws.Range("D34").Value = "word-manage-fields"
$ws.Range("E34").Value = "getFirstField"

$ws.Range("A35").Value = "FieldCollection"
$ws.Range("B35").Value = "getFirstOrNullObject"
$ws.Range("C35").Value = 1
$ws.Range("D35").Value = "word-manage-fields"
$ws.Range("E35").Value = "getFirstField"

$ws.Range("A36").Value = "FieldCollection"
$ws.Range("B36").Value = "items"
$ws.Range("C36").ClearContents()
$ws.Range("D36").Value = "word-manage-fields"
$ws.Range("E36").Value = "getAllFields"

# --- Normalise formatting ---------------------------------------------
# Rows 32-36 used to carry a couple of one-off cell styles (applyNumberFormat
# variants of the normal/centered styles) that only existed for this block.
# Re-apply the regular styles used by the rest of the table (copy format
# from existing rows) so those one-off styles become unused.
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A32:A36").PasteSpecial(-4122) | Out-Null

$ws.Range("D3").Copy() | Out-Null
$ws.Range("B32:E36").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# --- View state ----------------------------------------------------------
$ws.Range("E35").Select() | Out-Null
